# Apply the "Updated symbol list" edit to the cryptos sheet.
# Numeric-looking values are written with a leading apostrophe so that
# Excel keeps them as text (matching the original inlineStr cell type)
# instead of silently converting them to real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.98"
$ws.Range("G2").Value = "'6"
$ws.Range("G3").Value = "'6"
$ws.Range("D4").Value = "'5.207"
$ws.Range("G4").Value = "'6"
$ws.Range("D5").Value = "'0.05796"
$ws.Range("G5").Value = "'6"
$ws.Range("D6").Value = "'6.504"
$ws.Range("G6").Value = "'6"
$ws.Range("D7").Value = "'3.122"
$ws.Range("G7").Value = "'6"
$ws.Range("D8").Value = "'0.8156"
$ws.Range("G8").Value = "'6"
$ws.Range("D9").Value = "'0.8531"
$ws.Range("G9").Value = "'6"
$ws.Range("D10").Value = "'0.1358"
$ws.Range("G10").Value = "'6"
$ws.Range("D11").Value = "'0.06960"
$ws.Range("G11").Value = "'6"
$ws.Range("D12").Value = "'0.03168"
$ws.Range("G12").Value = "'6"
$ws.Range("D13").Value = "'0.02874"
$ws.Range("G13").Value = "'6"
$ws.Range("D14").Value = "'0.09381"
$ws.Range("G14").Value = "'6"
$ws.Range("D15").Value = "'3.741"
$ws.Range("G15").Value = "'6"
$ws.Range("D16").Value = "'0.001508"
$ws.Range("G16").Value = "'6"
$ws.Range("D17").Value = "'0.04677"
$ws.Range("G17").Value = "'6"
$ws.Range("D18").Value = "'0.0005961"
$ws.Range("G18").Value = "'6"
$ws.Range("D19").Value = "'0.006270"
$ws.Range("G19").Value = "'6"
$ws.Range("D20").Value = "'0.001236"
$ws.Range("G20").Value = "'6"
$ws.Range("D21").Value = "'0.004536"
$ws.Range("G21").Value = "'6"
$ws.Range("D22").Value = "'0.00006906"
$ws.Range("G22").Value = "'6"
$ws.Range("G23").Value = "'6"
$ws.Range("G24").Value = "'6"
$ws.Range("D25").Value = "'0.3173"
$ws.Range("G25").Value = "'6"
$ws.Range("G26").Value = "'6"
$ws.Range("D27").Value = "'0.1326"
$ws.Range("G27").Value = "'6"
$ws.Range("D28").Value = "'0.0002329"
$ws.Range("G28").Value = "'6"
$ws.Range("G29").Value = "'6"
$ws.Range("G30").Value = "'6"
$ws.Range("G31").Value = "'6"
$ws.Range("G32").Value = "'6"
$ws.Range("G33").Value = "'6"
$ws.Range("G34").Value = "'6"
$ws.Range("G35").Value = "'6"
$ws.Range("G36").Value = "'6"
$ws.Range("G37").Value = "'6"
$ws.Range("G38").Value = "'6"
$ws.Range("G39").Value = "'6"
$ws.Range("D40").Value = "'0.03656"
$ws.Range("G40").Value = "'6"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1052"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("G41").Value = "'6"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002852"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("G42").Value = "'6"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.002994"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("G43").Value = "'6"
$ws.Range("D44").Value = "'0.007479"
$ws.Range("G44").Value = "'6"
$ws.Range("D45").Value = "'0.00005279"
$ws.Range("G45").Value = "'6"
$ws.Range("G46").Value = "'6"
$ws.Range("D47").Value = "'0.3000"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("G47").Value = "'6"
$ws.Range("D48").Value = "'0.002340"
$ws.Range("G48").Value = "'6"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("G49").Value = "'6"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("G50").Value = "'6"
$ws.Range("G51").Value = "'6"
